# "test P7 with -10 percent"
# Rewrites the P5-scenario result values to the P7 (-10%) scenario values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "general": summary/objective values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Cells.Item(3, 2).Value  = 665.1125686385368    # objValue
$ws.Cells.Item(4, 2).Value  = 0.01800012588500977  # runtime
$ws.Cells.Item(6, 2).Value  = 50.56256919166783    # Z1
$ws.Cells.Item(7, 2).Value  = 0                    # Z2
$ws.Cells.Item(8, 2).Value  = 0                    # Z3
$ws.Cells.Item(9, 2).Value  = 440.0799994468689    # Z4
$ws.Cells.Item(10, 2).Value = 174.47               # Z5

# ---------------------------------------------------------------
# Sheet "x": column B (j) permuted
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("x")
$xVals = @(3,5,9,4,6,1,11,7,12,8,13,10,2)
for ($i = 0; $i -lt $xVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $xVals[$i]
}

# ---------------------------------------------------------------
# Sheet "U": column B (t) a few values changed
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Cells.Item(2, 2).Value  = 3
$ws.Cells.Item(4, 2).Value  = 2
$ws.Cells.Item(10, 2).Value = 2

# ---------------------------------------------------------------
# Sheet "TBar": column B (TBar) rows 3-15 updated (row 2 = 0 stays)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$tbarVals = @(
    35.43983387113133,
    34.33859424237852,
    27.8846700740982,
    41.48407906801076,
    30,
    38.68744466519142,
    33.99741166256366,
    38.86249131247691,
    28.59643311523818,
    30,
    45,
    42.75608190549734,
    44.07879268413863
)
for ($i = 0; $i -lt $tbarVals.Length; $i++) {
    $ws.Cells.Item($i + 3, 2).Value = $tbarVals[$i]
}

# ---------------------------------------------------------------
# Sheet "y": only header row remains, rows 2-11 removed
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("y")
$ws.Rows("2:11").Delete()

# ---------------------------------------------------------------
# Sheet "Q": column C (Q) values updated across many rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Q")
$qVals = @{
     7 = 150.2950000000002
     8 = 148.4950000000002
     9 = 140.8650000000002
    10 = 151.1300000000002
    11 = 149.8
    12 = 213.1800000000005
    13 = 212.6550000000005
    14 = 220.5800000000005
    15 = 220.0600000000005
    16 = 211.5400000000005
    17 = 40.35
    18 = 30.90499999999942
    19 = 27.59499999999942
    20 = 31.97499999999942
    21 = 33.99499999999942
    22 = 222.9249999999997
    23 = 228.55
    24 = 223.1849999999997
    25 = 235.0849999999997
    26 = 230.6149999999997
    29 = 72.41000000000022
    32 = 320.5600000000004
    33 = 324.8950000000003
    34 = 314.3700000000003
    35 = 328.8950000000004
    36 = 306.3200000000004
    37 = 152.9550000000011
    38 = 161.0400000000011
    39 = 150.3850000000011
    40 = 162.4150000000011
    41 = 154.7450000000011
    42 = 220.8800000000005
    43 = 240.2
    44 = 210.2050000000005
    45 = 226.1850000000005
    46 = 213.6750000000005
    47 = 73.75500000000051
    48 = 73.29500000000051
    49 = 74.13500000000052
    50 = 78.97000000000051
    51 = 71.2800000000005
    52 = 131.7049999999998
    53 = 131.9749999999998
    54 = 135.8649999999998
    55 = 139.6249999999998
    56 = 124.1949999999998
    57 = 222.9249999999997
    58 = 228.55
    59 = 223.1849999999997
    60 = 235.0849999999997
    61 = 230.6149999999997
    62 = 320.5600000000004
    63 = 324.8950000000003
    64 = 314.3700000000003
    65 = 328.8950000000004
    66 = 306.3200000000004
    67 = 220.8800000000005
    68 = 240.2
    69 = 210.2050000000005
    70 = 226.1850000000005
    71 = 213.6750000000005
}
foreach ($r in $qVals.Keys) {
    $ws.Cells.Item($r, 3).Value = $qVals[$r]
}

# ---------------------------------------------------------------
# Sheet "R": column C (R) rows 7-11 updated
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R")
$rVals = @{
     7 = 45.55999994277954
     8 = 49.89499998092651
     9 = 39.36999988555908
    10 = 53.89499998092651
    11 = 31.3199999332428
}
foreach ($r in $rVals.Keys) {
    $ws.Cells.Item($r, 3).Value = $rVals[$r]
}

# ---------------------------------------------------------------
# Sheet "L": column C (L) rows 2-6, 12-16, 42-46 updated
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$lVals = @{
     2 = 0
     3 = 0
     4 = 0
     5 = 0
     6 = 0
    12 = 6.57
    13 = 5.2
    14 = 7.32
    15 = 5.515
    16 = 5.44
    42 = 11.73
    43 = 14.67
    44 = 5.58
    45 = 11.76
    46 = 13.45
}
foreach ($r in $lVals.Keys) {
    $ws.Cells.Item($r, 3).Value = $lVals[$r]
}

# ---------------------------------------------------------------
# Sheet "rho": rows 7-11 removed, A2:A6 value 4 -> 6
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("rho")
$ws.Rows("7:11").Delete()
$ws.Range("A2:A6").Value = 6

# ---------------------------------------------------------------
# Sheet "alpha": rows 2-6 removed, only header remains
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("alpha")
$ws.Rows("2:6").Delete()
